$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 97 (shifts existing rows 97-117 down to 98-118)
$ws.Rows.Item(97).Insert()

# Populate the newly inserted row 97 with the new weekly data point
$ws.Cells.Item(97, 1).Value = 3
$ws.Cells.Item(97, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(97, 3).Value = "Coquimbo"
$ws.Cells.Item(97, 4).Value = 44551
$ws.Cells.Item(97, 5).Value = 5
$ws.Cells.Item(97, 6).Value = 100112026
$ws.Cells.Item(97, 7).Value = "Haba"
$ws.Cells.Item(97, 8).Value = "Sin especificar"
$ws.Cells.Item(97, 9).Value = "Primera"
$ws.Cells.Item(97, 10).Value = 105
$ws.Cells.Item(97, 11).Value = 8000
$ws.Cells.Item(97, 12).Value = 8500
$ws.Cells.Item(97, 13).Value = 8238
$ws.Cells.Item(97, 14).Value = "$/malla 25 kilos"
$ws.Cells.Item(97, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(97, 16).Value = 330
$ws.Cells.Item(97, 17).Value = 25
$ws.Cells.Item(97, 18).Value = "Hortaliza"

# Make sure the date cell keeps the same number format style as the rest of column D
$ws.Cells.Item(97, 4).NumberFormat = $ws.Cells.Item(98, 4).NumberFormat
